# "Copy with over ride Facitly"
#
# The Scenarios sheet's D9:G9 cells are text cells (stored as shared
# strings, numberformat General) that hold the sample values 235/23/48/11.
# The edit overwrites them - via a copy + paste-special-values - with the
# new sample values 888/777/999/666, keeping the cells as *text* (t="s")
# and keeping their existing cell style (s="12") untouched.
#
# A plain ".Value = '888'" assignment would make Excel auto-convert the
# numeric-looking string back into a Number, which is not what happened
# here (the cells stay text). Typing the values into a Text-formatted
# scratch cell and then doing Copy + PasteSpecial(xlPasteValues) onto the
# target reproduces the original "paste values, overriding the existing
# content" behaviour while preserving the destination's own formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# Use a far-away, unused part of the HELP sheet as scratch space for the
# copy source so nothing else on the workbook is permanently touched.
$scratch = $wb.Worksheets.Item("HELP")
$scratchRange = $scratch.Range("Z1:AC1")

$scratchRange.NumberFormat = "@"
$scratch.Range("Z1").Value = "888"
$scratch.Range("AA1").Value = "777"
$scratch.Range("AB1").Value = "999"
$scratch.Range("AC1").Value = "666"

$scratchRange.Copy()
$ws.Range("D9:G9").PasteSpecial(-4163)   # -4163 = xlPasteValues

$excel.CutCopyMode = $false
$scratchRange.Clear()

# Leave the selection where the user ended up after the paste.
$ws.Activate()
[void]$ws.Range("D13").Select()
